# Auto-generated PowerShell COM-interop script to apply the diff
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Define the new/changed text content as variables ----
$sEvalPartialHeader = @'
evaluator_partial_correctness
'@

$sLlmResponse10 = @'
To find the shortest path from node A to node J, we can use the breadth-first search (BFS) algorithm.
Starting from node A, we explore all its neighboring nodes (nodes B and C). Then, we explore the neighboring nodes of B and C, and continue this process until we reach node J.
The steps to find the shortest path are as follows:
1. Create an empty queue and enqueue node A.
2. Create an empty visited set to keep track of visited nodes.
3. Create an empty parent dictionary to store the parent node for each visited node.
4. While the queue is not empty, do the following:
   - Dequeue a node from the queue.
   - If the dequeued node is node J, we have found the shortest path.
   - Otherwise, for each neighboring node of the dequeued node that has not been visited, do the following:
     - Enqueue the neighboring node.
     - Mark the neighboring node as visited.
     - Set the parent of the neighboring node as the dequeued node in the parent dictionary.
5. If we have reached this step, it means that there is no path from node A to node J.
Using the adjacency matrix provided, we can implement the BFS algorithm to find the shortest path from node A to node J.
'@

$sWrong = @'
Wrong
'@

$sOutput010 = @'
Output: 0/10
'@

$sPrompt20 = @'
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node T?
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
    
'@

$sSolution20 = @'
A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T
'@

$sLlmResponse20 = @'
To find the shortest path from node A to node T, we can use Dijkstra's algorithm.
1. Initialize an array `dist` with size 20 to keep track of the shortest distance from node A to each other node. Start by setting the distance from A to A as 0 and the distance from A to all other nodes as infinity.
   dist = [0, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞, ∞]
2. Initialize a set `visited` to keep track of the nodes we have visited. Start by adding node A to the set.
   visited = {A}
3. While the set `visited` does not contain all nodes:
   - Find the node `u` in the set `visited` with the minimum distance from node A. This can be done by iterating through all nodes in the set and selecting the node with the minimum value in the `dist` array.
   - Add the node `u` to the set of visited nodes.
   - For each neighbor `v` of `u` that is not in the set `visited`:
     - Calculate the new distance `new_dist` from node A to `v` through `u` by adding the value in the adjacency matrix for the connection between `u` and `v` to the distance from A to `u`.
     - If `new_dist` is smaller than the current distance in the `dist` array for node `v`, update the `dist` array with the new distance.
4. After the algorithm terminates, the shortest distance from node A to node T will be the value in the `dist` array for node T.
In this case, the shortest path from node A to node T has a distance of 2.
Note: The actual path from node A to node T will depend on the actual connections between the nodes in the graph.
'@

$sOutput120 = @'
Output: 1/20
'@

$sLlmResponse20j = @'
To find the shortest path from node A to node T, we can use Dijkstra's algorithm. 
1. Start by setting the distance from node A to itself as 0, and all other distances as infinity. 
   - Distance[A] = 0
   - Distance[B] = infinity
   - Distance[C] = infinity
   ...
   - Distance[T] = infinity
2. Set the current node as A and mark it as visited. 
3. For each neighbor of the current node, calculate the distance from A to that neighbor through the current node. If this distance is shorter than the previously recorded distance for that neighbor, update it. 
   - For example, if current node is A and the neighbor is B: 
     - Distance[B] = min(Distance[B], Distance[A] + 1) = min(infinity, 0 + 1) = 1
   - Repeat this for all neighbors of the current node. 
4. Once all neighbors of the current node have been updated, mark the current node as visited. 
5. Select the unvisited node with the smallest distance as the new current node and repeat steps 3 to 5 until the destination node (T) is visited. 
6. The shortest path from A to T can be found by backtracking from T to A using the recorded distances. 
The step-by-step application of Dijkstra's algorithm to find the shortest path from A to T may be quite extensive to describe. Please let me know if a more detailed explanation is required.
'@

# ---- Step 1: Update sheet1 (o_10): add header E1 and set row2 values ----
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)
$ws1.Range("E1").Value = $sEvalPartialHeader

$ws1.Range("C2").Value = $sLlmResponse10
$ws1.Range("D2").Value = $sWrong
$ws1.Range("E2").Value = $sOutput010
$ws1.Rows.Item(2).AutoFit()

# ---- Step 2: Duplicate sheet1 -> o_20 (keeps identical header style/borders) ----
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "o_20"

$ws2.Range("A2").Value = $sPrompt20
$ws2.Range("B2").Value = $sSolution20
$ws2.Range("C2").Value = $sLlmResponse20
$ws2.Range("D2").Value = $sWrong
$ws2.Range("E2").Value = $sOutput120
$ws2.Rows.Item(2).AutoFit()

# ---- Step 3: Duplicate o_20 -> o_20_jumbled ----
$ws2.Copy([System.Reflection.Missing]::Value, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "o_20_jumbled"

$ws3.Range("C2").Value = $sLlmResponse20j
$ws3.Rows.Item(2).AutoFit()

# ---- Step 4: Restore the originally active sheet/tab ----
$ws1.Activate()
